# Apply the "Updated cryptos list" data refresh to the crypto table (rows 2-51).
# Column layout: A=rank index, B=coin name, C=link, D=price, E=volume(1h) change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.520.16'
$ws.Range('E2').Value = '  +4.08%  '
# Row 3
$ws.Range('D3').Value = '2.069.11'
$ws.Range('E3').Value = '  +3.55%  '
# Row 4
$ws.Range('E4').Value = '  +0.12%  '
# Row 5
$ws.Range('D5').Value = '''252.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.75%  '
# Row 6
$ws.Range('D6').Value = '''0.652'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.58%  '
# Row 7
$ws.Range('D7').Value = '''65.69'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +9.42%  '
# Row 8
$ws.Range('E8').Value = '  -0.02%  '
# Row 9
$ws.Range('E9').Value = '  +8.56%  '
# Row 10
$ws.Range('D10').Value = '''59.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.83%  '
# Row 11
$ws.Range('D11').Value = '''0.0827'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.00%  '
# Row 12
$ws.Range('E12').Value = '  +0.09%  '
# Row 13
$ws.Range('D13').Value = '''0.926'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.20%  '
# Row 14
$ws.Range('D14').Value = '''23.85'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +22.51%  '
# Row 15
$ws.Range('D15').Value = '''14.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.02%  '
# Row 16
$ws.Range('D16').Value = '2.373.85'
$ws.Range('E16').Value = '  +3.75%  '
# Row 17
$ws.Range('D17').Value = '''5.69'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.02%  '
# Row 18
$ws.Range('D18').Value = '2.059.54'
$ws.Range('E18').Value = '  +3.29%  '
# Row 19
$ws.Range('D19').Value = '37.436.24'
$ws.Range('E19').Value = '  +4.11%  '
# Row 20
$ws.Range('D20').Value = '''73.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.20%  '
# Row 21
$ws.Range('D21').Value = '0.0₃0909'
$ws.Range('E21').Value = '  +6.52%  '
# Row 22
$ws.Range('D22').Value = '''5.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +4.88%  '
# Row 23
$ws.Range('D23').Value = '''240.11'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.71%  '
# Row 24
$ws.Range('E24').Value = '  +0.07%  '
# Row 25
$ws.Range('D25').Value = '''2.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.32%  '
# Row 26
$ws.Range('D26').Value = '''2.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.60%  '
# Row 27
$ws.Range('D27').Value = '''10.08'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.89%  '
# Row 28
$ws.Range('D28').Value = '''20.83'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.19%  '
# Row 29
$ws.Range('E29').Value = '  -2.07%  '
# Row 30
$ws.Range('E30').Value = '  +28.97%  '
# Row 31
$ws.Range('E31').Value = '  +2.48%  '
# Row 32
$ws.Range('D32').Value = '''5.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.09%  '
# Row 33
$ws.Range('E33').Value = '  +5.64%  '
# Row 34
$ws.Range('D34').Value = '''0.0631'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.17%  '
# Row 35
$ws.Range('D35').Value = '''4.70'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.19%  '
# Row 36
$ws.Range('D36').Value = '''2.50'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.46%  '
# Row 37
$ws.Range('D37').Value = '''6.46'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.80%  '
# Row 38
$ws.Range('E38').Value = '  +0.20%  '
# Row 39
$ws.Range('E39').Value = '  +3.06%  '
# Row 40
$ws.Range('E40').Value = '  +31.02%  '
# Row 41
$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D41').Value = '''0.103'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.80%  '
# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.29'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.83%  '
# Row 43
$ws.Range('E43').Value = '  +4.61%  '
# Row 44
$ws.Range('E44').Value = '  +5.21%  '
# Row 45
$ws.Range('D45').Value = '''17.34'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.19%  '
# Row 46
$ws.Range('E46').Value = '  +2.20%  '
# Row 47
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').Value = '''8.00'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.78%  '
# Row 48
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''95.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.67%  '
# Row 49
$ws.Range('D49').Value = '1.401.98'
$ws.Range('E49').Value = '  +2.58%  '
# Row 50
$ws.Range('D50').Value = '''2.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.93%  '
# Row 51
$ws.Range('D51').Value = '''46.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.07%  '
